# arrumando correção de forma cobrança
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the corrected "forma cobrança" count for row 2 (ALL)
$ws.Range("B2").Value = 255

# Remove the test rows that were added by mistake (TESTE / Gamer / All family),
# clearing both the label (column A) and the value (column B) so the rows
# fall back to being just empty, styled placeholder cells like row 29.
$ws.Range("A17").ClearContents()
$ws.Range("B17").ClearContents()
$ws.Range("A18").ClearContents()
$ws.Range("B18").ClearContents()
$ws.Range("A19").ClearContents()
$ws.Range("B19").ClearContents()

# Restore the cursor/selection to C3
$ws.Range("C3").Select()
